$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P1")

$ws.Range("W3").Value = 1.341918521015387
$ws.Range("X3").Value = 1.39668576191086
$ws.Range("W4").Value = 1.471813954999379
$ws.Range("X4").Value = 1.472900366327911
$ws.Range("W5").Value = 1.365348369177526
$ws.Range("X5").Value = 1.356252746272316
$ws.Range("W6").Value = 0.7305066460837168
$ws.Range("X6").Value = 0.7462544375031713
$ws.Range("T7").Value = 0
$ws.Range("I9").Value = 0.5384615384615384
$ws.Range("J9").Value = 0
$ws.Range("S11").Value = 0.9886104783599089
$ws.Range("T11").Value = 0.4419134396355353
$ws.Range("J15").Value = 0.2519893899204244
$ws.Range("I19").Value = 0.4382022471910112
$ws.Range("J19").Value = 0.348314606741573
$ws.Range("H55").Value = 0.8503401360544192
$ws.Range("J55").Value = 0.1533742331288344
$ws.Range("J59").Value = 0.3333333333333334
$ws.Range("H63").Value = 1
$ws.Range("I63").Value = 1
$ws.Range("J63").Value = 0.06478873239436622
$ws.Range("H67").Value = 1
$ws.Range("I67").Value = 0.6279069767441861
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("H76").Value = 0.1428571428571429
$ws.Range("H79").Value = 1
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = 0.3895348837209303
$ws.Range("I83").Value = 0.3846153846153846
$ws.Range("J83").Value = 0.3076923076923077
$ws.Range("J101").Value = 0.01785714285714279
$ws.Range("T101").Value = 0.9821428571428572
$ws.Range("H111").Value = 1
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 0.2392156862745098
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("R112").Value = 1
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("H115").Value = 1
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = 0.1888888888888889
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("R116").Value = 1
$ws.Range("H119").Value = 1
$ws.Range("I119").Value = 1
$ws.Range("J119").Value = 0.3736263736263736
$ws.Range("I123").Value = 1
$ws.Range("J123").Value = 0.08571428571428574
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("R124").Value = 1
$ws.Range("S124").Value = 1
$ws.Range("T124").Value = 1
$ws.Range("H127").Value = 0.8579881656804733
$ws.Range("I127").Value = 0.8698224852071006
$ws.Range("J127").Value = 0.2280701754385965
$ws.Range("R127").Value = 1
$ws.Range("S127").Value = 1
$ws.Range("T127").Value = 1
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("H135").Value = 1
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = 0.3518005540166205
$ws.Range("R135").Value = 0.8423963133640553
$ws.Range("S135").Value = 0.8645161290322581
$ws.Range("T135").Value = 0.8866359447004608
$ws.Range("H136").Value = 1
$ws.Range("I136").Value = 1
$ws.Range("J136").Value = 0.2752043596730245
$ws.Range("R136").Value = 0.990990990990991
$ws.Range("H139").Value = 1
$ws.Range("I139").Value = 1
$ws.Range("J139").Value = 0.3594771241830066
$ws.Range("R139").Value = 0.828169014084507
$ws.Range("S139").Value = 0.8262910798122066
$ws.Range("T139").Value = 0.8413145539906103
$ws.Range("H140").Value = 1
$ws.Range("I140").Value = 1
$ws.Range("J140").Value = 0.3848580441640379
$ws.Range("R140").Value = 0.7076566125290024
$ws.Range("S140").Value = 0.777262180974478
$ws.Range("T140").Value = 0.8097447795823666
$ws.Range("J143").Value = 0.3625304136253041
$ws.Range("R143").Value = 0.8291316526610645
$ws.Range("S143").Value = 0.8141923436041083
$ws.Range("T143").Value = 0.830999066293184
$ws.Range("H144").Value = 1
$ws.Range("I144").Value = 1
$ws.Range("J144").Value = 0.4036939313984169
$ws.Range("R144").Value = 0.5868814729574223
$ws.Range("S144").Value = 0.6536248561565017
$ws.Range("T144").Value = 0.6881472957422324
$ws.Range("H147").Value = 1
$ws.Range("I147").Value = 1
$ws.Range("J147").Value = 0.3192771084337349
$ws.Range("S147").Value = 0.9663941871026339
$ws.Range("T147").Value = 0.9809264305177112
$ws.Range("H148").Value = 1
$ws.Range("I148").Value = 1
$ws.Range("J148").Value = 0.4399999999999999
$ws.Range("R148").Value = 0.7206385404789054
$ws.Range("S148").Value = 0.8164196123147093
$ws.Range("T148").Value = 0.8415051311288484
$ws.Range("H151").Value = 1
$ws.Range("I151").Value = 1
$ws.Range("J151").Value = 0.4167371090448013
$ws.Range("R151").Value = 0.7049873203719358
$ws.Range("S151").Value = 0.5562130177514792
$ws.Range("T151").Value = 0.5832628909551987
$ws.Range("H152").Value = 1
$ws.Range("I152").Value = 1
$ws.Range("J152").Value = 0.4449152542372882
$ws.Range("R152").Value = 0.8715415019762845
$ws.Range("S152").Value = 0.9387351778656127
$ws.Range("T152").Value = 0.9624505928853755
$ws.Range("H155").Value = 1
$ws.Range("I155").Value = 1
$ws.Range("J155").Value = 0.6833333333333333
$ws.Range("R155").Value = 0.4433333333333334
$ws.Range("S155").Value = 0.2916666666666667
$ws.Range("T155").Value = 0.3166666666666667
$ws.Range("H156").Value = 1
$ws.Range("I156").Value = 1
$ws.Range("J156").Value = 0.36
$ws.Range("H159").Value = 1
$ws.Range("I159").Value = 1
$ws.Range("J159").Value = 1
$ws.Range("R159").Value = 0.006546644844517169
$ws.Range("S159").Value = 0
$ws.Range("T159").Value = 0
$ws.Range("H160").Value = 1
$ws.Range("I160").Value = 1
$ws.Range("J160").Value = 0.8297872340425532
$ws.Range("R160").Value = 0.1814814814814815
$ws.Range("S160").Value = 0.1518518518518519
$ws.Range("T160").Value = 0.1851851851851851
$ws.Range("H161").Value = 1
$ws.Range("I161").Value = 1
$ws.Range("J161").Value = 1
$ws.Range("R161").Value = 0
$ws.Range("S161").Value = 0
$ws.Range("T161").Value = 0
$ws.Range("H163").Value = 1
$ws.Range("I163").Value = 1
$ws.Range("J163").Value = 1
$ws.Range("R163").Value = 0.2074829931972789
$ws.Range("S163").Value = 0
$ws.Range("T163").Value = 0
$ws.Range("H164").Value = 1
$ws.Range("I164").Value = 1
$ws.Range("J164").Value = 0.9318181818181818
$ws.Range("R164").Value = 0.4295302013422819
$ws.Range("S164").Value = 0.3911792905081496
$ws.Range("T164").Value = 0.4218600191754553
$ws.Range("H165").Value = 1
$ws.Range("I165").Value = 1
$ws.Range("J165").Value = 0.5114503816793894
$ws.Range("R165").Value = 0.366412213740458
$ws.Range("S165").Value = 0.366412213740458
$ws.Range("T165").Value = 0.4885496183206106
